$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$rows = @(2,3,4,5,6,10,13,14,18)
foreach ($r in $rows) {
    $ws.Range("M$r").ClearContents()
}

$ws.Range("M2").Select()
